$d = $word.ActiveDocument

# 1. Teléfono field: rename key to 'phone' and tweak placeholder/pattern
$rng = $d.Content
$found = $rng.Find.Execute(
    "[Teléfono;type='text';title='Teléfono';placeholder='+34 600 000 000';pattern='^\+?\d[\d\s\-]{7,}$';patternmsg='Formato de teléfono no válido';description='Admite prefijo internacional, espacios o guiones']",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0
)
if ($found) {
    $rng.Text = "[phone;type='text';title='Teléfono';placeholder='+34600000000';pattern='^[+]?[1-9][0-9]{1,14}$';patternmsg='Formato de teléfono no válido';description='Admite prefijo internacional, espacios o guiones']"
}

# 2. Cuerpo HTML field: rename key to 'body'
$rng = $d.Content
$found = $rng.Find.Execute(
    "[Cuerpo HTML;type='html';title='Cuerpo del documento (HTML)';description='Contenido enriquecido con formato']",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0
)
if ($found) {
    $rng.Text = "[body;type='html';title='Cuerpo del documento (HTML)';description='Contenido enriquecido con formato']"
}

# 3. Unidades field: rename key to 'units'
$rng = $d.Content
$found = $rng.Find.Execute(
    "[Unidades;type='number';title='Unidades';placeholder='0..20';minvalue='0';maxvalue='20';pattern='^(?:0|[1-9]|1[0-9]|20)$';patternmsg='Debe ser un entero entre 0 y 20']",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0
)
if ($found) {
    $rng.Text = "[units;type='number';title='Unidades';placeholder='0..20';minvalue='0';maxvalue='20';pattern='^(?:0|[1-9]|1[0-9]|20)$';patternmsg='Debe ser un entero entre 0 y 20']"
}

# 4. Last (empty) "Título 2" paragraph gets an explicit eastAsia font hint on its run properties
$lastPara = $d.Paragraphs.Last
$lastPara.Range.Font.NameFarEast = "eastAsia"
